# Update scripts wuth new tpm
# Applies new TPM-derived values to the LR-pairs sheet (Angpt4-Tek)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.816276414292558
$ws.Range("J2").Value = 0.816276414292558
$ws.Range("M2").Value = 71.05094633333333
$ws.Range("N2").Value = 213.152839
$ws.Range("O2").Value = 0.8240565632932695
$ws.Range("P2").Value = 0.8240565632932696
$ws.Range("Q2").Value = 46.08539638180956
$ws.Range("R2").Value = 414.768567436286
$ws.Range("S2").Value = 0.6726579366592784
$ws.Range("T2").Value = 0.6726579366592785

# Row 3
$ws.Range("I3").Value = 0.816276414292558
$ws.Range("J3").Value = 0.816276414292558
$ws.Range("O3").Value = 0.1323102827659759
$ws.Range("P3").Value = 0.132310282765976
$ws.Range("S3").Value = 0.1080017631902453
$ws.Range("T3").Value = 0.1080017631902453

# Row 4
$ws.Range("I4").Value = 0.816276414292558
$ws.Range("J4").Value = 0.816276414292558
$ws.Range("M4").Value = 3.762092333333333
$ws.Range("N4").Value = 11.286277
$ws.Range("O4").Value = 0.04363315394075456
$ws.Range("P4").Value = 0.04363315394075455
$ws.Range("Q4").Value = 2.440185885677555
$ws.Range("R4").Value = 21.961672971098
$ws.Range("S4").Value = 0.03561671444303433
$ws.Range("T4").Value = 0.03561671444303433

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1459893333333333
$ws.Range("H5").Value = 0.437968
$ws.Range("I5").Value = 0.183723585707442
$ws.Range("J5").Value = 0.183723585707442
$ws.Range("M5").Value = 71.05094633333333
$ws.Range("N5").Value = 213.152839
$ws.Range("O5").Value = 0.8240565632932695
$ws.Range("P5").Value = 0.8240565632932696
$ws.Range("Q5").Value = 10.37268028790578
$ws.Range("R5").Value = 93.35412259115201
$ws.Range("S5").Value = 0.1513986266339911
$ws.Range("T5").Value = 0.1513986266339911

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1459893333333333
$ws.Range("H6").Value = 0.437968
$ws.Range("I6").Value = 0.183723585707442
$ws.Range("J6").Value = 0.183723585707442
$ws.Range("O6").Value = 0.1323102827659759
$ws.Range("P6").Value = 0.132310282765976
$ws.Range("Q6").Value = 1.665434538193778
$ws.Range("R6").Value = 14.988910843744
$ws.Range("S6").Value = 0.02430851957573067
$ws.Range("T6").Value = 0.02430851957573067

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1459893333333333
$ws.Range("H7").Value = 0.437968
$ws.Range("I7").Value = 0.183723585707442
$ws.Range("J7").Value = 0.183723585707442
$ws.Range("M7").Value = 3.762092333333333
$ws.Range("N7").Value = 11.286277
$ws.Range("O7").Value = 0.04363315394075456
$ws.Range("P7").Value = 0.04363315394075455
$ws.Range("Q7").Value = 0.5492253516817777
$ws.Range("R7").Value = 4.943028165136
$ws.Range("S7").Value = 0.008016439497720232
$ws.Range("T7").Value = 0.00801643949772023
